$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.378.44"
$ws.Range("E2").Value = "  -1.10%  "
$ws.Range("D3").Value = "1.708.20"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "224.24"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.51%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5334"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.13%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.003"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2674"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.59%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06606"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.62%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.94"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.52%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07633"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.92%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.550"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.11%  "
$ws.Range("D13").Value = "1.702.51"
$ws.Range("E13").Value = "  -3.04%  "
$ws.Range("D14").Value = "1.945.79"
$ws.Range("E14").Value = "  -1.46%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5766"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Value = "0.0₅8179"
$ws.Range("E16").Value = "  -2.89%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "67.73"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.17%  "
$ws.Range("D18").Value = "27.349.03"
$ws.Range("E18").Value = "  -1.26%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "216.78"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -5.24%  "
$ws.Range("E20").Value = "  +0.04%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.668"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.45%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.44"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -4.22%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.963"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -4.23%  "
$ws.Range("E24").Value = "  -0.01%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "142.16"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -4.20%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.742"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.64%  "
$ws.Range("E27").Value = "  -2.79%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.259"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.83%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "16.26"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -5.04%  "
$ws.Range("E30").Value = "  -5.43%  "
$ws.Range("E31").Value = "  -1.80%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.491"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -5.69%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.424"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.73%  "
$ws.Range("E34").Value = "  -2.98%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.871"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.37%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9488"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.81%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.414"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.07%  "
$ws.Range("E38").Value = "  -2.09%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01633"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.34%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.849"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.82%  "
$ws.Range("D41").Value = "1.044.61"
$ws.Range("E41").Value = "  -0.52%  "
$ws.Range("E42").Value = "  -0.05%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8383"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.44%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "100.91"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.68%  "
$ws.Range("D45").Value = "1.851.79"
$ws.Range("E45").Value = "  -1.48%  "
$ws.Range("D46").Value = "0.0₈118"
$ws.Range("E46").Value = "  +2.06%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "57.92"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.84%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4512"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.69%  "
$ws.Range("E49").Value = "  +0.39%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.083"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.41%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05236"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.90%  "
